$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "Andrea Conzatti "
$ws.Range("B11").Value = "Gabriel Melis | Demobusters"
$ws.Range("C11").Value = "Carlo Stedile | MAI UNA GIOIA"
$ws.Range("D11").Value = "Davide Deimichei | SdrumALA"
$ws.Range("E11").Value = "Tommaso Caliari | RSA United"
$ws.Range("F11").Value = "FEDERICO NICOLODI | U.S. Guarna"
